$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LogisticRegression - Obesity")
$ws2 = $wb.Worksheets.Item("Summary")

# -----------------------------------------------------------------
# Sheet2 "Summary": insert two new rows at the top. Row 1 becomes a
# new title row, row 2 stays blank, and all previously existing rows
# (old 1-17) shift down by two (new 3-19).
# -----------------------------------------------------------------
$ws2.Range("A1:D2").EntireRow.Insert()

# New title in A1, using the same look as the other section titles
# (bold red font, no fill/border) that already exists in the workbook
# on sheet1!A2 (style index 3).
$ws1.Range("A2").Copy()
$ws2.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("A1").Value = "Using LR, one hot encoding and ngram(1,2)"

# Fill in the previously empty Accuracy(%) column (D) for the single
# -feature data rows (now rows 4-11; D4 already had a value of 95.2).
$ws2.Range("D5").Value = 94.8
$ws2.Range("D6").Value = 86.7
$ws2.Range("D7").Value = 89.7
$ws2.Range("D8").Value = 94
$ws2.Range("D9").Value = 93.6
$ws2.Range("D10").Value = 89.7
$ws2.Range("D11").Value = 93.7

# Fill in the Accuracy(%) column for the cumulative-feature rows
# (now rows 13-19).
$ws2.Range("D13").Value = 95
$ws2.Range("D14").Value = 94
$ws2.Range("D15").Value = 94
$ws2.Range("D16").Value = 93.7
$ws2.Range("D17").Value = 94.1
$ws2.Range("D18").Value = 93.8
$ws2.Range("D19").Value = 94.1

# Widen column A to fit the new, longer disease/combo labels.
$ws2.Columns.Item(1).ColumnWidth = 38.65

# -----------------------------------------------------------------
# View-state updates (selected cell on each sheet + which sheet tab
# is active). Restore Summary as the active sheet at the end so the
# saved tabSelected flag stays on Summary, matching the original file.
# -----------------------------------------------------------------
$ws1.Range("I11").Select()
$ws2.Range("D13").Select()
$ws2.Activate()
